# Apply the "Added DarkData to Cell Result export, Added Suffix to CellGroups" edit.
#
# Summary of changes:
#  - templateResultSheet (sheet 1): add a "Cells" / "#cells" column (N),
#    move the #UIChartLightDark marker from E6 to F6, add Dark Voltage/DarkCurrent
#    header+key columns (C9:D9, C10:D10), and add a "Light Data"/"Dark Data"
#    label row (32) above the relocated #UIChart/#UIChartDark markers (F33/M33).
#  - TemplateCellGroup (sheet 2): add a "Cells" / "#cells" column (N).
#  - CellGroupsTemplate (sheet 3): add a "CELLS_MEAN" / "#cells" column (AI).
#  - Column widths on sheet1/sheet2 collapse to a single uniform width.
#  - The active sheet/tab moves back to templateResultSheet.

$wb = $excel.ActiveWorkbook

$wsResult = $wb.Worksheets.Item("templateResultSheet")
$wsCellGroup = $wb.Worksheets.Item("TemplateCellGroup")
$wsCellGroupsTemplate = $wb.Worksheets.Item("CellGroupsTemplate")

# ---------------------------------------------------------------------------
# templateResultSheet
# ---------------------------------------------------------------------------

# New "Cells" column (header row 1 + key row 2)
$wsResult.Range("N1").Value = "Cells"
$wsResult.Range("N2").Value = "#cells"

# #UIChartLightDark marker moves from E6 to F6
$wsResult.Range("E6").ClearContents()
$wsResult.Range("F6").Value = "#UIChartLightDark"

# New Dark Voltage / DarkCurrent header + key columns
$wsResult.Range("C9").Value = "Dark Voltage[V]"
$wsResult.Range("D9").Value = "DarkCurrent[A]"
$wsResult.Range("C10").Value = "#darkVoltageData"
$wsResult.Range("D10").Value = "#darkCurrentData"

# New row 32 with Light Data / Dark Data labels
$wsResult.Range("F32").Value = "Light Data"
$wsResult.Range("M32").Value = "Dark Data"

# Row 33 markers move from D33/H33 to F33/M33
$wsResult.Range("D33").ClearContents()
$wsResult.Range("H33").ClearContents()
$wsResult.Range("F33").Value = "#UIChart"
$wsResult.Range("M33").Value = "#UIChartDark"

# Columns A:N collapse to a single uniform width
$wsResult.Range("A1:N1").EntireColumn.ColumnWidth = 14.13

# ---------------------------------------------------------------------------
# TemplateCellGroup
# ---------------------------------------------------------------------------

$wsCellGroup.Range("N1").Value = "Cells"
$wsCellGroup.Range("N2").Value = "#cells"

# Columns B:N collapse to a single uniform width (A keeps its own width)
$wsCellGroup.Range("B1:N1").EntireColumn.ColumnWidth = 14.13

# ---------------------------------------------------------------------------
# CellGroupsTemplate
# ---------------------------------------------------------------------------

$wsCellGroupsTemplate.Range("AI1").Value = "CELLS_MEAN"
$wsCellGroupsTemplate.Range("AI2").Value = "#cells"

# ---------------------------------------------------------------------------
# Active sheet moves back to templateResultSheet
# ---------------------------------------------------------------------------

$wsResult.Activate()
